$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Update the reporting period dates (row 8)
$ws.Range("B8").Value = 44743
$ws.Range("C8").Value = 44926

# Update "last update" / "validation" dates (row 8)
$ws.Range("AG8").Value = 44936
$ws.Range("AH8").Value = 44936

# Update the selected cell in the active sheet view
$ws.Range("C11").Select()
